$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text values look like plain numbers need to be forced to
# Text format first, otherwise Excel auto-converts them to numeric values
# and strips significant trailing/leading zeros (e.g. "0.100" -> 0.1).
$textCells = @("D5","D6","D8","D10","D11","D12","D13","D14","D16","D17","D21","D22","D23","D24","D25","D28","D29","D30","D31","D32","D33","D34","D36","D37","D38","D40","D41","D46","D47","D48","D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "70.699.26"
$ws.Range("E2").Value = "  +4.43%  "
$ws.Range("D3").Value = "3.553.03"
$ws.Range("E3").Value = "  +3.54%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "603.28"
$ws.Range("E5").Value = "  +4.51%  "
$ws.Range("D6").Value = "172.93"
$ws.Range("E6").Value = "  +5.19%  "
$ws.Range("D7").Value = "3.544.84"
$ws.Range("E7").Value = "  +3.52%  "
$ws.Range("D8").Value = "0.613"
$ws.Range("E8").Value = "  +2.27%  "
$ws.Range("E9").Value = "  -0.07%  "
$ws.Range("D10").Value = "0.193"
$ws.Range("E10").Value = "  +5.90%  "
$ws.Range("D11").Value = "7.34"
$ws.Range("E11").Value = "  +9.71%  "
$ws.Range("D12").Value = "0.586"
$ws.Range("E12").Value = "  +3.64%  "
$ws.Range("D13").Value = "46.31"
$ws.Range("E13").Value = "  +0.86%  "
$ws.Range("D14").Value = "0.0000276"
$ws.Range("E14").Value = "  +3.02%  "
$ws.Range("D15").Value = "4.130.21"
$ws.Range("E15").Value = "  +3.56%  "
$ws.Range("D16").Value = "8.31"
$ws.Range("E16").Value = "  +1.20%  "
$ws.Range("D17").Value = "606.47"
$ws.Range("E17").Value = "  -0.42%  "
$ws.Range("D18").Value = "3.556.27"
$ws.Range("E18").Value = "  +3.61%  "
$ws.Range("D19").Value = "70.741.16"
$ws.Range("E19").Value = "  +4.34%  "
$ws.Range("E20").Value = "  +1.55%  "
$ws.Range("D21").Value = "17.31"
$ws.Range("D22").Value = "0.877"
$ws.Range("E22").Value = "  +1.52%  "
$ws.Range("D23").Value = "9.27"
$ws.Range("E23").Value = "  -14.64%  "
$ws.Range("D24").Value = "15.68"
$ws.Range("E24").Value = "  +2.37%  "
$ws.Range("D25").Value = "96.42"
$ws.Range("E25").Value = "  +1.86%  "
$ws.Range("E26").Value = "  +0.87%  "
$ws.Range("E27").Value = "  -0.07%  "
$ws.Range("D28").Value = "2.60"
$ws.Range("E28").Value = "  +2.06%  "
$ws.Range("D29").Value = "33.98"
$ws.Range("E29").Value = "  +6.44%  "
$ws.Range("D30").Value = "9.05"
$ws.Range("E30").Value = "  +2.16%  "
$ws.Range("D31").Value = "721.54"
$ws.Range("E31").Value = "  +20.59%  "
$ws.Range("D32").Value = "3.04"
$ws.Range("E32").Value = "  +1.16%  "
$ws.Range("D33").Value = "8.21"
$ws.Range("E33").Value = "  -0.65%  "
$ws.Range("D34").Value = "7.02"
$ws.Range("E34").Value = "  +4.58%  "
$ws.Range("E35").Value = "  +1.38%  "
$ws.Range("D36").Value = "0.100"
$ws.Range("E36").Value = "  +0.94%  "
$ws.Range("D37").Value = "3.57"
$ws.Range("E37").Value = "  +6.26%  "
$ws.Range("D38").Value = "10.74"
$ws.Range("E38").Value = "  +1.60%  "
$ws.Range("E39").Value = "  +10.97%  "
$ws.Range("D40").Value = "56.91"
$ws.Range("E40").Value = "  +0.82%  "
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  +0.05%  "
$ws.Range("E42").Value = "  +7.38%  "
$ws.Range("D43").Value = "3.363.96"
$ws.Range("E43").Value = "  +0.59%  "
$ws.Range("E44").Value = "  -0.12%  "
$ws.Range("E45").Value = "  +3.65%  "
$ws.Range("D46").Value = "32.51"
$ws.Range("E46").Value = "  +1.42%  "
$ws.Range("D47").Value = "2.92"
$ws.Range("E47").Value = "  +8.82%  "
$ws.Range("D48").Value = "2.58"
$ws.Range("E48").Value = "  +5.07%  "
$ws.Range("E49").Value = "  +2.47%  "
$ws.Range("D50").Value = "133.92"
$ws.Range("E50").Value = "  +1.22%  "
$ws.Range("E51").Value = "  -0.02%  "
